$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

$ws.Cells.Item($row, 1).Value = 36
$ws.Cells.Item($row, 2).Value = "11:37"
$ws.Cells.Item($row, 3).Value = "riya-morankar"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "edit1 to main"

# Column F holds a date-like text string ("2025-06-18") that must stay as
# plain text (matching the other rows), not get auto-converted to a real
# Excel date serial number. Force text format, assign, then clear the
# formatting override again so no stray style is left behind on the cell.
$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "2025-06-18"
$ws.Cells.Item($row, 6).ClearFormats()
